# Daily update (8-10 data): append new demographic rows dated 2020-08-10 (serial 44053)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("age")
$ws.Cells.Item(886, 1).Value = 44053
$ws.Cells.Item(886, 2).Value = "<1 year"
$ws.Cells.Item(886, 3).Value = 143
$ws.Cells.Item(886, 4).Value = 0.003503871410369499
$ws.Cells.Item(886, 5).Value = 0
$ws.Cells.Item(886, 6).Value = 2
$ws.Cells.Item(886, 7).Value = 0.0002355712603062426
$ws.Cells.Item(886, 8).Value = 0
$ws.Cells.Item(887, 1).Value = 44053
$ws.Cells.Item(887, 2).Value = "1-9 years"
$ws.Cells.Item(887, 3).Value = 651
$ws.Cells.Item(887, 4).Value = 0.01595119082622758
$ws.Cells.Item(887, 5).Value = 1
$ws.Cells.Item(887, 6).Value = 5
$ws.Cells.Item(887, 7).Value = 0.0005889281507656066
$ws.Cells.Item(887, 8).Value = 0
$ws.Cells.Item(888, 1).Value = 44053
$ws.Cells.Item(888, 2).Value = "10-19 years"
$ws.Cells.Item(888, 3).Value = 1682
$ws.Cells.Item(888, 4).Value = 0.04121336861707341
$ws.Cells.Item(888, 5).Value = 8
$ws.Cells.Item(888, 6).Value = 9
$ws.Cells.Item(888, 7).Value = 0.001060070671378092
$ws.Cells.Item(888, 8).Value = 1
$ws.Cells.Item(889, 1).Value = 44053
$ws.Cells.Item(889, 2).Value = "20-29 years"
$ws.Cells.Item(889, 3).Value = 7567
$ws.Cells.Item(889, 4).Value = 0.1854111535822797
$ws.Cells.Item(889, 5).Value = 61
$ws.Cells.Item(889, 6).Value = 63
$ws.Cells.Item(889, 7).Value = 0.007420494699646643
$ws.Cells.Item(889, 8).Value = 0
$ws.Cells.Item(890, 1).Value = 44053
$ws.Cells.Item(890, 2).Value = "30-39 years"
$ws.Cells.Item(890, 3).Value = 8252
$ws.Cells.Item(890, 4).Value = 0.2021954327158679
$ws.Cells.Item(890, 5).Value = 74
$ws.Cells.Item(890, 6).Value = 175
$ws.Cells.Item(890, 7).Value = 0.02061248527679623
$ws.Cells.Item(890, 8).Value = 0
$ws.Cells.Item(891, 1).Value = 44053
$ws.Cells.Item(891, 2).Value = "40-49 years"
$ws.Cells.Item(891, 3).Value = 7674
$ws.Cells.Item(891, 4).Value = 0.188032931490738
$ws.Cells.Item(891, 5).Value = 67
$ws.Cells.Item(891, 6).Value = 459
$ws.Cells.Item(891, 7).Value = 0.05406360424028268
$ws.Cells.Item(891, 8).Value = 0
$ws.Cells.Item(892, 1).Value = 44053
$ws.Cells.Item(892, 2).Value = "50-59 years"
$ws.Cells.Item(892, 3).Value = 6992
$ws.Cells.Item(892, 4).Value = 0.1713221601489758
$ws.Cells.Item(892, 5).Value = 44
$ws.Cells.Item(892, 6).Value = 981
$ws.Cells.Item(892, 7).Value = 0.115547703180212
$ws.Cells.Item(892, 8).Value = 4
$ws.Cells.Item(893, 1).Value = 44053
$ws.Cells.Item(893, 2).Value = "60-64 years"
$ws.Cells.Item(893, 3).Value = 2711
$ws.Cells.Item(893, 4).Value = 0.06642654121336862
$ws.Cells.Item(893, 5).Value = 13
$ws.Cells.Item(893, 6).Value = 799
$ws.Cells.Item(893, 7).Value = 0.09411071849234394
$ws.Cells.Item(893, 8).Value = 1
$ws.Cells.Item(894, 1).Value = 44053
$ws.Cells.Item(894, 2).Value = "65-69 years"
$ws.Cells.Item(894, 3).Value = 1928
$ws.Cells.Item(894, 4).Value = 0.04724100754679996
$ws.Cells.Item(894, 5).Value = 7
$ws.Cells.Item(894, 6).Value = 972
$ws.Cells.Item(894, 7).Value = 0.1144876325088339
$ws.Cells.Item(894, 8).Value = 1
$ws.Cells.Item(895, 1).Value = 44053
$ws.Cells.Item(895, 2).Value = "70-74 years"
$ws.Cells.Item(895, 3).Value = 1146
$ws.Cells.Item(895, 4).Value = 0.02807997647750662
$ws.Cells.Item(895, 5).Value = 1
$ws.Cells.Item(895, 6).Value = 1101
$ws.Cells.Item(895, 7).Value = 0.1296819787985866
$ws.Cells.Item(895, 8).Value = 4
$ws.Cells.Item(896, 1).Value = 44053
$ws.Cells.Item(896, 2).Value = "75-79 years"
$ws.Cells.Item(896, 3).Value = 764
$ws.Cells.Item(896, 4).Value = 0.01871998431833774
$ws.Cells.Item(896, 5).Value = 0
$ws.Cells.Item(896, 6).Value = 1054
$ws.Cells.Item(896, 7).Value = 0.1241460541813899
$ws.Cells.Item(896, 8).Value = 5
$ws.Cells.Item(897, 1).Value = 44053
$ws.Cells.Item(897, 2).Value = "80+ years"
$ws.Cells.Item(897, 3).Value = 1271
$ws.Cells.Item(897, 4).Value = 0.03114280113692051
$ws.Cells.Item(897, 5).Value = 1
$ws.Cells.Item(897, 6).Value = 2870
$ws.Cells.Item(897, 7).Value = 0.3380447585394582
$ws.Cells.Item(897, 8).Value = 15
$ws.Cells.Item(898, 1).Value = 44053
$ws.Cells.Item(898, 2).Value = "Unknown"
$ws.Cells.Item(898, 3).Value = 31
$ws.Cells.Item(898, 4).Value = 0.0007595805155346467
$ws.Cells.Item(898, 5).Value = 0
$ws.Cells.Item(898, 6).Value = 0
$ws.Cells.Item(898, 7).Value = 0
$ws.Cells.Item(898, 8).Value = 0

$ws = $wb.Worksheets.Item("gender")
$ws.Cells.Item(206, 1).Value = 44053
$ws.Cells.Item(206, 2).Value = "Female"
$ws.Cells.Item(206, 3).Value = 15063
$ws.Cells.Item(206, 4).Value = 0.3690826227580124
$ws.Cells.Item(206, 5).Value = 39
$ws.Cells.Item(206, 6).Value = 3456
$ws.Cells.Item(206, 7).Value = 0.4070671378091873
$ws.Cells.Item(206, 8).Value = 14
$ws.Cells.Item(207, 1).Value = 44053
$ws.Cells.Item(207, 2).Value = "Male"
$ws.Cells.Item(207, 3).Value = 24853
$ws.Cells.Item(207, 4).Value = 0.6089630500833089
$ws.Cells.Item(207, 5).Value = 237
$ws.Cells.Item(207, 6).Value = 5034
$ws.Cells.Item(207, 7).Value = 0.5929328621908128
$ws.Cells.Item(207, 8).Value = 17
$ws.Cells.Item(208, 1).Value = 44053
$ws.Cells.Item(208, 2).Value = "Unknown"
$ws.Cells.Item(208, 3).Value = 896
$ws.Cells.Item(208, 4).Value = 0.02195432715867882
$ws.Cells.Item(208, 5).Value = 1
$ws.Cells.Item(208, 6).Value = 0
$ws.Cells.Item(208, 7).Value = 0
$ws.Cells.Item(208, 8).Value = 0

$ws = $wb.Worksheets.Item("race")
$ws.Cells.Item(404, 1).Value = 44053
$ws.Cells.Item(404, 2).Value = "Asian"
$ws.Cells.Item(404, 3).Value = 723
$ws.Cells.Item(404, 4).Value = 0.01771537783004999
$ws.Cells.Item(404, 5).Value = 1
$ws.Cells.Item(404, 6).Value = 190
$ws.Cells.Item(404, 7).Value = 0.02237926972909305
$ws.Cells.Item(404, 8).Value = 0
$ws.Cells.Item(405, 1).Value = 44053
$ws.Cells.Item(405, 2).Value = "Black"
$ws.Cells.Item(405, 3).Value = 6179
$ws.Cells.Item(405, 4).Value = 0.1514015485641478
$ws.Cells.Item(405, 5).Value = 73
$ws.Cells.Item(405, 6).Value = 1086
$ws.Cells.Item(405, 7).Value = 0.1279151943462898
$ws.Cells.Item(405, 8).Value = 0
$ws.Cells.Item(406, 1).Value = 44053
$ws.Cells.Item(406, 2).Value = "Hispanic"
$ws.Cells.Item(406, 3).Value = 16150
$ws.Cells.Item(406, 4).Value = 0.3957169459962756
$ws.Cells.Item(406, 5).Value = 115
$ws.Cells.Item(406, 6).Value = 4405
$ws.Cells.Item(406, 7).Value = 0.5188457008244994
$ws.Cells.Item(406, 8).Value = 19
$ws.Cells.Item(407, 1).Value = 44053
$ws.Cells.Item(407, 2).Value = "Other"
$ws.Cells.Item(407, 3).Value = 84
$ws.Cells.Item(407, 4).Value = 0.002058218171126139
$ws.Cells.Item(407, 5).Value = -93
$ws.Cells.Item(407, 6).Value = 58
$ws.Cells.Item(407, 7).Value = 0.006831566548881037
$ws.Cells.Item(407, 8).Value = 0
$ws.Cells.Item(408, 1).Value = 44053
$ws.Cells.Item(408, 2).Value = "Unknown"
$ws.Cells.Item(408, 3).Value = 6542
$ws.Cells.Item(408, 4).Value = 0.1602959913750858
$ws.Cells.Item(408, 5).Value = 103
$ws.Cells.Item(408, 6).Value = 4
$ws.Cells.Item(408, 7).Value = 0.0004711425206124853
$ws.Cells.Item(408, 8).Value = 0
$ws.Cells.Item(409, 1).Value = 44053
$ws.Cells.Item(409, 2).Value = "White"
$ws.Cells.Item(409, 3).Value = 11134
$ws.Cells.Item(409, 4).Value = 0.2728119180633147
$ws.Cells.Item(409, 5).Value = 78
$ws.Cells.Item(409, 6).Value = 2747
$ws.Cells.Item(409, 7).Value = 0.3235571260306243
$ws.Cells.Item(409, 8).Value = 12
